$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rightmost values to append (after the shift) for rows 2-6.
# Rows 7-16 simply lose their rightmost value (no replacement).
$newValues = @{
    2 = -0.001489938197266189
    3 = -0.182031752916177
    4 = -0.2180070093596886
    5 = 0.4767206611340558
    6 = 1.178844253737389
}

for ($row = 2; $row -le 16; $row++) {
    # Determine the last populated column in this row (B=2 .. K=11)
    $lastCol = 1
    for ($col = 2; $col -le 11; $col++) {
        if ($ws.Cells.Item($row, $col).Value2 -ne $null) {
            $lastCol = $col
        }
    }

    # Shift every value one column to the left, dropping column B's value.
    for ($col = 2; $col -le $lastCol; $col++) {
        if ($col -lt $lastCol) {
            $ws.Cells.Item($row, $col).Value2 = $ws.Cells.Item($row, $col + 1).Value2
        } else {
            if ($newValues.ContainsKey($row)) {
                $ws.Cells.Item($row, $col).Value2 = $newValues[$row]
            } else {
                $ws.Cells.Item($row, $col).Value2 = $null
            }
        }
    }
}
